# Automatische test-sync: 2025-08-03 14:47:50
# Appends the new "Testmail #7" row to the Logs sheet and bumps the
# "Inkoop / Bestellingen" tally on the Dashboard sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = 17
$logs.Cells.Item($newRow, 1).Value = "Is dit artikel nog op voorraad?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #7: Is dit artikel nog op voorraad?"
$logs.Cells.Item($newRow, 4).Value = "Inkoop / Bestellingen"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-03 14:46:50"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(5, 2).Value = 2

# The new row extends every existing conditional-formatting band by one row
# (…2:…16 -> …2:…17) while leaving the rules themselves untouched.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "16")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "17")
    foreach ($fc in $oldRange.FormatConditions) {
        $fc.ModifyAppliesToRange($newRange)
    }
}

